$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp label in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 00:50"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 186046
$ws.Range("C4").Value = 22258
$ws.Range("E4").Value = 175892
$ws.Range("G4").Value = 666
$ws.Range("H4").Value = 3807

# Row 21 - Israel
$ws.Range("B21").Value = 4831
$ws.Range("C21").Value = 136
$ws.Range("E21").Value = 4587

# Row 22 - Australia
$ws.Range("B22").Value = 4712
$ws.Range("C22").Value = 252
$ws.Range("E22").Value = 4355

# Row 25 - Chequia
$ws.Range("B25").Value = 3308
$ws.Range("C25").Value = 307
$ws.Range("E25").Value = 3232

# Row 135 - Guatemala
$ws.Range("D135").Value = 12
$ws.Range("E135").Value = 25
